$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 235 is a brand-new row; populate the "constant" columns by copying row 234,
# matching every other data row (A,B,C,E,F,G,H,I,N,O,Q,R are identical across all rows).
$ws.Cells.Item(235, 1).Value2 = $ws.Cells.Item(234, 1).Value2
$ws.Cells.Item(235, 2).Value2 = $ws.Cells.Item(234, 2).Value2
$ws.Cells.Item(235, 3).Value2 = $ws.Cells.Item(234, 3).Value2
$ws.Cells.Item(235, 5).Value2 = $ws.Cells.Item(234, 5).Value2
$ws.Cells.Item(235, 6).Value2 = $ws.Cells.Item(234, 6).Value2
$ws.Cells.Item(235, 7).Value2 = $ws.Cells.Item(234, 7).Value2
$ws.Cells.Item(235, 8).Value2 = $ws.Cells.Item(234, 8).Value2
$ws.Cells.Item(235, 9).Value2 = $ws.Cells.Item(234, 9).Value2
$ws.Cells.Item(235, 14).Value2 = $ws.Cells.Item(234, 14).Value2
$ws.Cells.Item(235, 15).Value2 = $ws.Cells.Item(234, 15).Value2
$ws.Cells.Item(235, 17).Value2 = $ws.Cells.Item(234, 17).Value2
$ws.Cells.Item(235, 18).Value2 = $ws.Cells.Item(234, 18).Value2
# Fecha (column D) carries a date NumberFormat on every row; match it on the new row too.
$ws.Cells.Item(235, 4).NumberFormat = $ws.Cells.Item(234, 4).NumberFormat

# A new weekly price observation was inserted at the top of the series (row 39);
# every subsequent Fecha/Volumen/Precio tuple (D,J,K,L,M,P) shifts down by one row,
# appending a brand-new trailing row (235) that holds the former last row's (234) data.
$shiftedRows = @(
    @{Row=39; D=44558; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=40; D=44497; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=41; D=44435; J=810; K=1500; L=1500; M=1500; P=500}
    @{Row=42; D=44445; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=43; D=44496; J=150; K=1500; L=1500; M=1500; P=500}
    @{Row=44; D=44249; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=45; D=44295; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=46; D=44356; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=47; D=44438; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=48; D=44244; J=110; K=1500; L=1500; M=1500; P=500}
    @{Row=49; D=44342; J=260; K=1500; L=1500; M=1500; P=500}
    @{Row=50; D=44202; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=51; D=44509; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=52; D=44448; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=53; D=44333; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=54; D=44159; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=55; D=44336; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=56; D=44515; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=57; D=44343; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=58; D=44245; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=59; D=44397; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=60; D=44523; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=61; D=44524; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=62; D=44216; J=80; K=1500; L=1500; M=1500; P=500}
    @{Row=63; D=44369; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=64; D=44410; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=65; D=44319; J=190; K=1500; L=1500; M=1500; P=500}
    @{Row=66; D=44186; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=67; D=44189; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=68; D=44355; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=69; D=44351; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=70; D=44442; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=71; D=44363; J=130; K=1500; L=1500; M=1500; P=500}
    @{Row=72; D=44406; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=73; D=44215; J=130; K=1500; L=1500; M=1500; P=500}
    @{Row=74; D=44392; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=75; D=44517; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=76; D=44544; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=77; D=44284; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=78; D=44484; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=79; D=44300; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=80; D=44426; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=81; D=44201; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=82; D=44358; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=83; D=44382; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=84; D=44431; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=85; D=44530; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=86; D=44179; J=48; K=2000; L=2000; M=2000; P=667}
    @{Row=87; D=44477; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=88; D=44512; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=89; D=44557; J=80; K=1500; L=1500; M=1500; P=500}
    @{Row=90; D=44203; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=91; D=44384; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=92; D=44236; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=93; D=44487; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=94; D=44320; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=95; D=44294; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=96; D=44326; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=97; D=44407; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=98; D=44315; J=130; K=1500; L=1500; M=1500; P=500}
    @{Row=99; D=44483; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=100; D=44505; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=101; D=44274; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=102; D=44348; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=103; D=44309; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=104; D=44488; J=150; K=1500; L=1500; M=1500; P=500}
    @{Row=105; D=44214; J=110; K=1500; L=1500; M=1500; P=500}
    @{Row=106; D=44172; J=110; K=1500; L=1500; M=1500; P=500}
    @{Row=107; D=44546; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=108; D=44237; J=130; K=1500; L=1500; M=1500; P=500}
    @{Row=109; D=44252; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=110; D=44383; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=111; D=44312; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=112; D=44162; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=113; D=44349; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=114; D=44376; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=115; D=44554; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=116; D=44299; J=130; K=1500; L=1500; M=1500; P=500}
    @{Row=117; D=44246; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=118; D=44327; J=190; K=1500; L=1500; M=1500; P=500}
    @{Row=119; D=44316; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=120; D=44174; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=121; D=44277; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=122; D=44181; J=90; K=1500; L=1500; M=1500; P=500}
    @{Row=123; D=44482; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=124; D=44273; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=125; D=44364; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=126; D=44469; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=127; D=44265; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=128; D=44330; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=129; D=44280; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=130; D=44323; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=131; D=44209; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=132; D=44266; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=133; D=44403; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=134; D=44165; J=68; K=2000; L=2000; M=2000; P=667}
    @{Row=135; D=44267; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=136; D=44306; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=137; D=44263; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=138; D=44516; J=150; K=1500; L=1500; M=1500; P=500}
    @{Row=139; D=44427; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=140; D=44533; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=141; D=44495; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=142; D=44176; J=80; K=1500; L=1500; M=1500; P=500}
    @{Row=143; D=44539; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=144; D=44370; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=145; D=44475; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=146; D=44508; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=147; D=44291; J=89; K=1800; L=1800; M=1800; P=600}
    @{Row=148; D=44468; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=149; D=44379; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=150; D=44532; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=151; D=44498; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=152; D=44489; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=153; D=44210; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=154; D=44526; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=155; D=44271; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=156; D=44425; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=157; D=44476; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=158; D=44218; J=130; K=1500; L=1500; M=1500; P=500}
    @{Row=159; D=44250; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=160; D=44168; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=161; D=44447; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=162; D=44553; J=150; K=1500; L=1500; M=1500; P=500}
    @{Row=163; D=44167; J=150; K=1500; L=1500; M=1500; P=500}
    @{Row=164; D=44161; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=165; D=44328; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=166; D=44285; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=167; D=44418; J=150; K=1500; L=1500; M=1500; P=500}
    @{Row=168; D=44160; J=230; K=1500; L=1500; M=1500; P=500}
    @{Row=169; D=44434; J=140; K=1500; L=1500; M=1500; P=500}
    @{Row=170; D=44467; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=171; D=44231; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=172; D=44490; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=173; D=44259; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=174; D=44251; J=80; K=1500; L=1500; M=1500; P=500}
    @{Row=175; D=44341; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=176; D=44286; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=177; D=44279; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=178; D=44208; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=179; D=44264; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=180; D=44322; J=130; K=1500; L=1500; M=1500; P=500}
    @{Row=181; D=44491; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=182; D=44391; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=183; D=44396; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=184; D=44510; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=185; D=44232; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=186; D=44551; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=187; D=44386; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=188; D=44519; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=189; D=44420; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=190; D=44414; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=191; D=44543; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=192; D=44321; J=130; K=1500; L=1500; M=1500; P=500}
    @{Row=193; D=44385; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=194; D=44278; J=130; K=1500; L=1500; M=1500; P=500}
    @{Row=195; D=44308; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=196; D=44281; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=197; D=44474; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=198; D=44446; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=199; D=44350; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=200; D=44529; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=201; D=44405; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=202; D=44413; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=203; D=44238; J=130; K=1500; L=1500; M=1500; P=500}
    @{Row=204; D=44257; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=205; D=44411; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=206; D=44175; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=207; D=44196; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=208; D=44200; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=209; D=44459; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=210; D=44188; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=211; D=44258; J=230; K=1500; L=1500; M=1500; P=500}
    @{Row=212; D=44298; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=213; D=44432; J=150; K=1500; L=1500; M=1500; P=500}
    @{Row=214; D=44428; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=215; D=44329; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=216; D=44452; J=190; K=1500; L=1500; M=1500; P=500}
    @{Row=217; D=44270; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=218; D=44195; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=219; D=44473; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=220; D=44398; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=221; D=44302; J=130; K=1500; L=1500; M=1500; P=500}
    @{Row=222; D=44511; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=223; D=44239; J=120; K=1500; L=1500; M=1500; P=500}
    @{Row=224; D=44344; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=225; D=44461; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=226; D=44463; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=227; D=44357; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=228; D=44371; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=229; D=44365; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=230; D=44454; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=231; D=44194; J=80; K=1500; L=1500; M=1500; P=500}
    @{Row=232; D=44313; J=130; K=1500; L=1500; M=1500; P=500}
    @{Row=233; D=44518; J=160; K=1500; L=1500; M=1500; P=500}
    @{Row=234; D=44540; J=180; K=1500; L=1500; M=1500; P=500}
    @{Row=235; D=44272; J=160; K=1500; L=1500; M=1500; P=500}
)

foreach ($item in $shiftedRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value2 = $item.D
    $ws.Cells.Item($r, 10).Value2 = $item.J
    $ws.Cells.Item($r, 11).Value2 = $item.K
    $ws.Cells.Item($r, 12).Value2 = $item.L
    $ws.Cells.Item($r, 13).Value2 = $item.M
    $ws.Cells.Item($r, 16).Value2 = $item.P
}
